$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/17/2024  Through  6/23/2024"

# --- Row 14 ---
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -82.352941176470

# --- Row 15 ---
$ws.Range("M15").Value = 0

# --- Row 16 ---
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 88
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 11.392405063291
$ws.Range("L16").Value = 6.024096385542
$ws.Range("M16").Value = -11.111111111111
$ws.Range("N16").Value = -76.902887139107

# --- Row 17 ---
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 119
$ws.Range("K17").Value = 11.764705882352
$ws.Range("L17").Value = 26.666666666666
$ws.Range("M17").Value = 43.010752688172
$ws.Range("N17").Value = -59.202453987730

# --- Row 18 ---
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 46
$ws.Range("K18").Value = -24.590163934426
$ws.Range("L18").Value = -46.511627906976
$ws.Range("M18").Value = -17.857142857142
$ws.Range("N18").Value = -89.327146171693

# --- Row 19 ---
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 243
$ws.Range("J19").Value = 167
$ws.Range("K19").Value = 45.508982035928
$ws.Range("L19").Value = 39.655172413793
$ws.Range("M19").Value = 68.75
$ws.Range("N19").Value = 23.979591836734

# --- Row 20 ---
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 31
$ws.Range("K20").Value = -12.903225806451
$ws.Range("L20").Value = 17.391304347826
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -60.869565217391

# --- Row 21 ---
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = 2.531645569620
$ws.Range("I21").Value = 546
$ws.Range("J21").Value = 463
$ws.Range("K21").Value = 17.92656587473
$ws.Range("L21").Value = 14.465408805031
$ws.Range("M21").Value = 33.170731707317
$ws.Range("N21").Value = -61.977715877437

# --- Row 22 (G22/H22 become suppressed "0"/"***.*" placeholders like sibling cells) ---
$ws.Range("G22").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

# --- Row 23 ---
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 45
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 12.5
$ws.Range("L23").Value = -2.173913043478
$ws.Range("M23").Value = 221.428571428571

# --- Row 24 ---
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -23.333333333333
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -3.921568627450
$ws.Range("I24").Value = 595
$ws.Range("J24").Value = 604
$ws.Range("K24").Value = -1.490066225165
$ws.Range("L24").Value = 9.778597785977
$ws.Range("M24").Value = 32.8125

# --- Row 25 ---
$ws.Range("C25").Value = 13
$ws.Range("E25").Value = -18.75
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = 1.694915254237
$ws.Range("I25").Value = 371
$ws.Range("J25").Value = 364
$ws.Range("K25").Value = 1.923076923076
$ws.Range("L25").Value = 31.560283687943

# --- Row 26 ---
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 30
$ws.Range("I26").Value = 196
$ws.Range("J26").Value = 199
$ws.Range("K26").Value = -1.507537688442
$ws.Range("L26").Value = -4.854368932038
$ws.Range("M26").Value = -20.967741935483

# --- Row 28 (D28/E28 become suppressed "0"/"***.*" placeholders like sibling cells) ---
$ws.Range("D28").Value = "'0"
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = 10.526315789473
$ws.Range("L28").Value = -25

# --- Row 29 ---
$ws.Range("L29").Value = -22.222222222222
$ws.Range("N29").Value = -86

# --- Row 30 ---
$ws.Range("L30").Value = -33.333333333333
$ws.Range("N30").Value = -86.666666666666
